# Adds two new "weekly report" blocks (rows 112-120 and 122-130) to the
# bottom of the sheet, mirroring the existing block structure (rows 102-110)
# for formatting/merges, then fills in the new text content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Clone formatting + merges from the last existing block (102:110)
#    down onto the two new blocks. Copy also recreates the A*:D* merge
#    on the title/summary rows automatically.
# ---------------------------------------------------------------------
$ws.Range("A102:D110").Copy($ws.Range("A112:D120"))
$ws.Range("A102:D110").Copy($ws.Range("A122:D130"))

# ---------------------------------------------------------------------
# 2) Row heights (auto-fit differs because the new text differs from the
#    copied source text, so set them explicitly to match the authored
#    workbook).
# ---------------------------------------------------------------------
$ws.Rows.Item(112).RowHeight = 22.5
$ws.Rows.Item(113).RowHeight = 22.5
$ws.Rows.Item(114).RowHeight = 68
$ws.Rows.Item(115).RowHeight = 22.5
$ws.Rows.Item(116).RowHeight = 22.5
$ws.Rows.Item(117).RowHeight = 22.5
$ws.Rows.Item(118).RowHeight = 22.5
$ws.Rows.Item(119).RowHeight = 22.5
$ws.Rows.Item(120).RowHeight = 22.5
$ws.Rows.Item(121).RowHeight = 29
$ws.Rows.Item(122).RowHeight = 22.5
$ws.Rows.Item(123).RowHeight = 22.5
$ws.Rows.Item(124).RowHeight = 45
$ws.Rows.Item(125).RowHeight = 22.5
$ws.Rows.Item(126).RowHeight = 22.5
$ws.Rows.Item(127).RowHeight = 22.5
$ws.Rows.Item(128).RowHeight = 45
$ws.Rows.Item(129).RowHeight = 22.5
$ws.Rows.Item(130).RowHeight = 22.5

# ---------------------------------------------------------------------
# 3) Fix the two cells whose style diverges from the 102:110 template
#    (completion-status cells that need wrapText in the new block).
# ---------------------------------------------------------------------
$ws.Range("C115").Copy($ws.Range("C116"))
$ws.Range("C115").Copy($ws.Range("C119"))

# ---------------------------------------------------------------------
# 4) Block 1 (日期：2017.10.11 第七周 周三) — rows 112-120
# ---------------------------------------------------------------------
$ws.Range("A112").Value = "日期：2017.10.11 第七周 周三"

$ws.Range("A113").Value = "人员"
$ws.Range("B113").Value = "计划任务"
$ws.Range("C113").Value = "完成情况"
$ws.Range("D113").Value = "备注"

$ws.Range("A114").Value = "李杰"
$ws.Range("B114").Value = "继续编写后台“用户管理”的数据交互，并完成数据库接口文档"
$ws.Range("C114").Value = "未完成"

$ws.Range("A115").Value = "周振朋"
$ws.Range("B115").Value = "尝试开发“首页”界面"
$ws.Range("C115").Value = "未完成"

$ws.Range("A116").Value = "禤锦辉"
$ws.Range("B116").Value = "尝试开发“买卖”界面"
$ws.Range("C116").Value = "未完成"

$ws.Range("A117").Value = "柯新钿"
$ws.Range("B117").Value = "尝试开发“登录”界面"
$ws.Range("C117").Value = "完成了一部分"

$ws.Range("A118").Value = "冯文雄"
$ws.Range("B118").Value = "完成web接口设计文档"
$ws.Range("C118").Value = "未完成"

$ws.Range("A119").Value = "阿卜力孜"
$ws.Range("B119").Value = "尝试开发“我的”界面"
$ws.Range("C119").Value = "未完成"

$ws.Range("A120").Value = "总结：未恢复状态，完成情况很差"

# ---------------------------------------------------------------------
# 5) Block 2 (日期：2017.10.16 第八周 周一) — rows 122-130
# ---------------------------------------------------------------------
$ws.Range("A122").Value = "日期：2017.10.16 第八周 周一"

$ws.Range("A123").Value = "人员"
$ws.Range("B123").Value = "计划任务"
$ws.Range("C123").Value = "完成情况"
$ws.Range("D123").Value = "备注"

$ws.Range("A124").Value = "李杰"
$ws.Range("B124").Value = "完善数据库接口文档和web接口文档"
$ws.Range("C124").ClearContents()

$ws.Range("A125").Value = "周振朋"
$ws.Range("B125").Value = "继续开发“首页”界面"
$ws.Range("C125").ClearContents()

$ws.Range("A126").Value = "禤锦辉"
$ws.Range("B126").Value = "继续开发“买卖”界面"
$ws.Range("C126").ClearContents()

$ws.Range("A127").Value = "柯新钿"
$ws.Range("B127").Value = "开发“注册”界面"
$ws.Range("C127").ClearContents()

$ws.Range("A128").Value = "冯文雄"
$ws.Range("B128").Value = "继续编写后台“商品类别管理”的数据交互"
$ws.Range("C128").ClearContents()

$ws.Range("A129").Value = "阿卜力孜"
$ws.Range("B129").Value = "继续开发“我的”界面"
$ws.Range("C129").ClearContents()

$ws.Range("A130").Value = "总结："

# ---------------------------------------------------------------------
# 6) Selection / view bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("A130:D130").Select()
